# Fruta / hortaliza, semanal
# The data rows (2-30) get reshuffled to a new row order (weekly refresh of the
# consolidated subset). No new values are introduced; every resulting row is an
# exact copy of one of the original rows, just placed at a different row number.
# Column A is constant (always 1) for every row, so we only need to move B..R.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row number -> source row number (values taken from source row
# in the *original* sheet and placed into target row in the final sheet).
$rowMap = @{
    2  = 7
    3  = 20
    4  = 9
    5  = 13
    6  = 2
    7  = 5
    8  = 6
    9  = 15
    10 = 28
    11 = 19
    12 = 18
    13 = 23
    14 = 24
    15 = 17
    16 = 27
    17 = 21
    18 = 29
    19 = 30
    20 = 16
    21 = 26
    22 = 10
    23 = 11
    24 = 12
    25 = 14
    26 = 22
    27 = 8
    28 = 3
    29 = 4
    30 = 25
}

$firstCol = 2   # column B
$lastCol  = 18  # column R

# Snapshot every source cell's value first (Value2 gives back plain scalars,
# e.g. date serials instead of COM wrapper objects) so that overwriting rows
# while iterating never clobbers data we still need to read later.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write each target row using the snapshotted source row data.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($targetRow, $c).Value2 = $sourceVals[$c]
    }
}
